$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed coin stats (price + 1h volume change), scraped on a later run of
# the GitHub Action. Row 47/48 also swapped rank (ordi <-> Aave) this time.
#
# Every cell below is plain text in the sheet (coinranking prices use "." as a
# thousands separator, e.g. "43.127.45", and keep trailing zeros, e.g. "5.60").
# "AsText" marks the few new values that *look* like a plain decimal number -
# left to Excel's auto-detection, `.Value = "5.60"` would silently store a
# Number (5.6, wrong type AND wrong text). For those cells we force
# `NumberFormat = "@"` before the write and restore the "Normal" style right
# after, so no stray number-format style is left behind on the cell.
$edits = @(
    @{ Cell = "D2"; Value = "43.127.45"; AsText = $false }
    @{ Cell = "E2"; Value = "  +0.94%  "; AsText = $false }
    @{ Cell = "D3"; Value = "2.376.02"; AsText = $false }
    @{ Cell = "E3"; Value = "  +0.16%  "; AsText = $false }
    @{ Cell = "E4"; Value = "  +0.03%  "; AsText = $false }
    @{ Cell = "D5"; Value = "316.75"; AsText = $true }
    @{ Cell = "E5"; Value = "  -3.24%  "; AsText = $false }
    @{ Cell = "D6"; Value = "109.33"; AsText = $true }
    @{ Cell = "E6"; Value = "  +10.67%  "; AsText = $false }
    @{ Cell = "E7"; Value = "  +0.54%  "; AsText = $false }
    @{ Cell = "E8"; Value = "  -0.04%  "; AsText = $false }
    @{ Cell = "D9"; Value = "0.623"; AsText = $true }
    @{ Cell = "E9"; Value = "  +0.61%  "; AsText = $false }
    @{ Cell = "D10"; Value = "41.24"; AsText = $true }
    @{ Cell = "E10"; Value = "  +3.76%  "; AsText = $false }
    @{ Cell = "E11"; Value = "  +1.28%  "; AsText = $false }
    @{ Cell = "D12"; Value = "8.65"; AsText = $true }
    @{ Cell = "E12"; Value = "  +3.21%  "; AsText = $false }
    @{ Cell = "D13"; Value = "1.01"; AsText = $true }
    @{ Cell = "E13"; Value = "  +0.76%  "; AsText = $false }
    @{ Cell = "E15"; Value = "  -1.38%  "; AsText = $false }
    @{ Cell = "D16"; Value = "2.734.33"; AsText = $false }
    @{ Cell = "E16"; Value = "  +0.09%  "; AsText = $false }
    @{ Cell = "D17"; Value = "2.404.71"; AsText = $false }
    @{ Cell = "E17"; Value = "  +1.31%  "; AsText = $false }
    @{ Cell = "D18"; Value = "43.122.60"; AsText = $false }
    @{ Cell = "E18"; Value = "  +0.92%  "; AsText = $false }
    @{ Cell = "D19"; Value = "7.74"; AsText = $true }
    @{ Cell = "E19"; Value = "  -0.88%  "; AsText = $false }
    @{ Cell = "E20"; Value = "  +1.35%  "; AsText = $false }
    @{ Cell = "D21"; Value = "76.58"; AsText = $true }
    @{ Cell = "E21"; Value = "  +2.03%  "; AsText = $false }
    @{ Cell = "E22"; Value = "  -2.66%  "; AsText = $false }
    @{ Cell = "D23"; Value = "269.60"; AsText = $true }
    @{ Cell = "E23"; Value = "  -0.60%  "; AsText = $false }
    @{ Cell = "D24"; Value = "2.35"; AsText = $true }
    @{ Cell = "E24"; Value = "  +1.62%  "; AsText = $false }
    @{ Cell = "E25"; Value = "  -0.60%  "; AsText = $false }
    @{ Cell = "E26"; Value = "  +0.38%  "; AsText = $false }
    @{ Cell = "D27"; Value = "11.52"; AsText = $true }
    @{ Cell = "E27"; Value = "  +0.91%  "; AsText = $false }
    @{ Cell = "D28"; Value = "23.65"; AsText = $true }
    @{ Cell = "E28"; Value = "  +0.29%  "; AsText = $false }
    @{ Cell = "D29"; Value = "2.26"; AsText = $true }
    @{ Cell = "E29"; Value = "  +2.13%  "; AsText = $false }
    @{ Cell = "D30"; Value = "37.37"; AsText = $true }
    @{ Cell = "E30"; Value = "  +6.67%  "; AsText = $false }
    @{ Cell = "D31"; Value = "169.06"; AsText = $true }
    @{ Cell = "E31"; Value = "  -2.08%  "; AsText = $false }
    @{ Cell = "E32"; Value = "  +1.59%  "; AsText = $false }
    @{ Cell = "D33"; Value = "6.19"; AsText = $true }
    @{ Cell = "E33"; Value = "  +5.63%  "; AsText = $false }
    @{ Cell = "D34"; Value = "2.94"; AsText = $true }
    @{ Cell = "E34"; Value = "  -5.71%  "; AsText = $false }
    @{ Cell = "E35"; Value = "  +15.97%  "; AsText = $false }
    @{ Cell = "E36"; Value = "  +1.11%  "; AsText = $false }
    @{ Cell = "D37"; Value = "4.74"; AsText = $true }
    @{ Cell = "E37"; Value = "  +3.74%  "; AsText = $false }
    @{ Cell = "D38"; Value = "0.0364"; AsText = $true }
    @{ Cell = "E38"; Value = "  +2.21%  "; AsText = $false }
    @{ Cell = "D39"; Value = "3.88"; AsText = $true }
    @{ Cell = "E39"; Value = "  +0.60%  "; AsText = $false }
    @{ Cell = "E40"; Value = "  -4.85%  "; AsText = $false }
    @{ Cell = "D41"; Value = "105.17"; AsText = $true }
    @{ Cell = "E41"; Value = "  +9.67%  "; AsText = $false }
    @{ Cell = "E42"; Value = "  +1.60%  "; AsText = $false }
    @{ Cell = "D43"; Value = "0.239"; AsText = $true }
    @{ Cell = "E43"; Value = "  +5.76%  "; AsText = $false }
    @{ Cell = "D44"; Value = "72.01"; AsText = $true }
    @{ Cell = "E44"; Value = "  +5.60%  "; AsText = $false }
    @{ Cell = "D45"; Value = "12.81"; AsText = $true }
    @{ Cell = "E45"; Value = "  +9.07%  "; AsText = $false }
    @{ Cell = "E46"; Value = "  +0.00%  "; AsText = $false }
    @{ Cell = "B47"; Value = "Aave"; AsText = $false }
    @{ Cell = "C47"; Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"; AsText = $false }
    @{ Cell = "D47"; Value = "114.96"; AsText = $true }
    @{ Cell = "E47"; Value = "  -0.77%  "; AsText = $false }
    @{ Cell = "B48"; Value = "ordi"; AsText = $false }
    @{ Cell = "C48"; Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"; AsText = $false }
    @{ Cell = "D48"; Value = "80.50"; AsText = $true }
    @{ Cell = "E48"; Value = "  +19.67%  "; AsText = $false }
    @{ Cell = "D49"; Value = "5.60"; AsText = $true }
    @{ Cell = "E49"; Value = "  +3.71%  "; AsText = $false }
    @{ Cell = "D50"; Value = "9.22"; AsText = $true }
    @{ Cell = "E50"; Value = "  +3.39%  "; AsText = $false }
    @{ Cell = "E51"; Value = "  +4.07%  "; AsText = $false }
)

foreach ($edit in $edits) {
    $rng = $ws.Range($edit.Cell)
    if ($edit.AsText) {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $edit.Value
    if ($edit.AsText) {
        $rng.Style = "Normal"
    }
}
